# TimeSheet_Week8 update — "Attendance up to date : Oct 11, 2021"
# Record 1 hour worked on Monday (column B) for the "Sponsor Meeting" row
# (row 9), then move the active selection to the Week Total cell for that
# column (B13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 1

$ws.Range("B13").Select()
